$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values that look numeric (e.g. "521.46").
# If assigned directly, Excel auto-converts them into real floating point
# numbers (losing the exact decimal text and trailing zeros, e.g. "1.00" -> 1).
# Force each target cell to Text format before assigning the string, then
# restore the default "Normal" style so the cell keeps no explicit style -
# matching the original file where these are plain text values.
$dRows = @(2,3,5,6,9,11,12,14,15,16,17,18,19,20,21,23,24,26,27,28,29,31,33,34,35,36,37,38,39,40,43,45,46,47,49,50,51)
foreach ($r in $dRows) { $ws.Range("D$r").NumberFormat = "@" }

$ws.Range("D2").Value = '58.792.29'
$ws.Range("D3").Value = '3.089.41'
$ws.Range("D5").Value = '521.46'
$ws.Range("D6").Value = '144.05'
$ws.Range("D9").Value = '7.36'
$ws.Range("D11").Value = '0.384'
$ws.Range("D12").Value = '3.619.67'
$ws.Range("D14").Value = '26.71'
$ws.Range("D15").Value = '0.0000167'
$ws.Range("D16").Value = '58.793.00'
$ws.Range("D17").Value = '3.087.96'
$ws.Range("D18").Value = '6.15'
$ws.Range("D19").Value = '12.92'
$ws.Range("D20").Value = '8.13'
$ws.Range("D21").Value = '344.35'
$ws.Range("D23").Value = '0.506'
$ws.Range("D24").Value = '65.75'
$ws.Range("D26").Value = '1.00'
$ws.Range("D27").Value = '0.0₃0924'
$ws.Range("D28").Value = '6.63'
$ws.Range("D29").Value = '7.27'
$ws.Range("D31").Value = '21.00'
$ws.Range("D33").Value = '154.40'
$ws.Range("D34").Value = '4.62'
$ws.Range("D35").Value = '6.13'
$ws.Range("D36").Value = '26.75'
$ws.Range("D37").Value = '1.30'
$ws.Range("D38").Value = '0.0686'
$ws.Range("D39").Value = '3.128.90'
$ws.Range("D40").Value = '3.92'
$ws.Range("D43").Value = '0.665'
$ws.Range("D45").Value = '2.287.48'
$ws.Range("D46").Value = '0.0255'
$ws.Range("D47").Value = '20.67'
$ws.Range("D49").Value = '6.00'
$ws.Range("D50").Value = '0.748'
$ws.Range("D51").Value = '264.52'

foreach ($r in $dRows) { $ws.Range("D$r").Style = "Normal" }

# "Volume(1h)" column (E) values contain "%" and surrounding spaces, so Excel
# keeps them as plain text without any extra coercion.
$ws.Range("E2").Value = '  +2.01%  '
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E5").Value = '  +0.99%  '
$ws.Range("E6").Value = '  +0.76%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +0.58%  '
$ws.Range("E9").Value = '  +0.31%  '
$ws.Range("E10").Value = '  +0.54%  '
$ws.Range("E11").Value = '  +2.72%  '
$ws.Range("E12").Value = '  +0.55%  '
$ws.Range("E13").Value = '  +0.95%  '
$ws.Range("E14").Value = '  +3.56%  '
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("E16").Value = '  +1.87%  '
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("E18").Value = '  -0.11%  '
$ws.Range("E19").Value = '  -1.23%  '
$ws.Range("E20").Value = '  -0.91%  '
$ws.Range("E21").Value = '  +1.87%  '
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("E23").Value = '  +0.89%  '
$ws.Range("E25").Value = '  -0.52%  '
$ws.Range("E26").Value = '  -0.31%  '
$ws.Range("E27").Value = '  -0.94%  '
$ws.Range("E28").Value = '  +2.50%  '
$ws.Range("E29").Value = '  +2.59%  '
$ws.Range("E30").Value = '  +1.29%  '
$ws.Range("E31").Value = '  +0.64%  '
$ws.Range("E32").Value = '  +1.90%  '
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("E34").Value = '  +2.05%  '
$ws.Range("E35").Value = '  +3.52%  '
$ws.Range("E36").Value = '  +0.24%  '
$ws.Range("E37").Value = '  +3.86%  '
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("E39").Value = '  +0.36%  '
$ws.Range("E40").Value = '  +1.17%  '
$ws.Range("E41").Value = '  -0.69%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("E43").Value = '  -1.08%  '
$ws.Range("E44").Value = '  +4.28%  '
$ws.Range("E45").Value = '  +0.38%  '
$ws.Range("E46").Value = '  +0.50%  '
$ws.Range("E47").Value = '  +1.59%  '
$ws.Range("E48").Value = '  +0.24%  '
$ws.Range("E49").Value = '  +2.01%  '
$ws.Range("E50").Value = '  +8.39%  '
$ws.Range("E51").Value = '  +11.42%  '
